$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "readability" row (row 145) entirely, shifting all rows below it up by one.
$ws.Rows.Item(145).Delete()
